# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2310"
#   "<header>_new" -> "<header>_FV2404"
# Also freeze the header row and wrap the data range in an Excel Table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (A1:J1 = "_old" -> "_FV2310", L1:U1 = "_new" -> "_FV2404"; K1 "diff" stays) ---
$ws.Range("A1").Value2 = "Segmentname_FV2310"
$ws.Range("B1").Value2 = "Segmentgruppe_FV2310"
$ws.Range("C1").Value2 = "Segment_FV2310"
$ws.Range("D1").Value2 = "Datenelement_FV2310"
$ws.Range("E1").Value2 = "Segment ID_FV2310"
$ws.Range("F1").Value2 = "Code_FV2310"
$ws.Range("G1").Value2 = "Qualifier_FV2310"
$ws.Range("H1").Value2 = "Beschreibung_FV2310"
$ws.Range("I1").Value2 = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value2 = "Bedingung_FV2310"

$ws.Range("L1").Value2 = "Segmentname_FV2404"
$ws.Range("M1").Value2 = "Segmentgruppe_FV2404"
$ws.Range("N1").Value2 = "Segment_FV2404"
$ws.Range("O1").Value2 = "Datenelement_FV2404"
$ws.Range("P1").Value2 = "Segment ID_FV2404"
$ws.Range("Q1").Value2 = "Code_FV2404"
$ws.Range("R1").Value2 = "Qualifier_FV2404"
$ws.Range("S1").Value2 = "Beschreibung_FV2404"
$ws.Range("T1").Value2 = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value2 = "Bedingung_FV2404"

# --- 2. Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the A1:U79 range into an Excel Table named "Table1" ---
$range = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
